$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-18: update date serial 45172 -> 45175
# (2023-09-03 -> 2023-09-06) for every data row.
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
